$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# D-column price cells are stored as literal text in the source data (they use
# "." as a thousands separator, e.g. "63.772.34"), so each write forces Text
# number-format before the assignment (otherwise values that are also valid
# numbers, like "568.86", get silently coerced to numeric -- dropping trailing
# zeros / switching to scientific notation for tiny values) and then restores
# the original (unstyled/General) look by copying the neighbouring Link cell's
# style, so no stray style index is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.669.42"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("E2").Value = "  -1.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.411.19"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.86"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.54"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.413.42"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = "  -0.38%  "

$ws.Range("E9").Value = "  -6.49%  "

$ws.Range("E10").Value = "  +0.86%  "

$ws.Range("E11").Value = "  -2.20%  "

$ws.Range("E12").Value = "  -3.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.008.90"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.11"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = "  -2.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = "  -7.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.750.56"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.428.30"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = "  -2.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = "  -3.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.58"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.57"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = "  +1.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.75"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "  -2.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.10"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  -1.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.518"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = "  -5.51%  "

$ws.Range("E26").Value = "  -3.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.66"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = "  -3.62%  "

$ws.Range("E28").Value = "  +0.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.992"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = "  -0.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.08"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = "  -0.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.39"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = "  -6.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.98"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = "  -1.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.95"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  -0.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.98"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = "  -3.11%  "

$ws.Range("E35").Value = "  -4.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.19"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = "  +0.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.832"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = "  +8.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = "  -2.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.15"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = "  -2.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.805.97"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = "  -1.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0722"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  -4.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.80"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.37"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = "  -6.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.40"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = "  -4.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.66"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("E46").Value = "  -2.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "326.04"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").Value = "  +3.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "  +8.18%  "

$ws.Range("E49").Value = "  -4.02%  "

$ws.Range("E50").Value = "  -4.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.30"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = "  -3.83%  "
